$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update two numeric "Level of Effort" values
$ws.Range("F11").Value = 20
$ws.Range("F16").Value = 40

# Insert a new blank row before row 28 (shifts rows 28-37 down to 29-38)
$ws.Rows("28:28").Insert()

# The inserted row is blank, but should carry the same look as the row
# above it (row 27) rather than Excel's bare default - copy formatting over.
$ws.Range("A27:H27").Copy()
$ws.Range("A28:H28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The rows that used to be 28, 29, 30 are now 29, 30, 31.
# Change their "Release" (column A) value from the number 2 to the text "?"
$ws.Range("A29").Value = "?"
$ws.Range("A30").Value = "?"
$ws.Range("A31").Value = "?"

# Restore the selection to match the post-edit state
$ws.Range("H18").Select()
